{"js": "// Replace the answer text in every cell of the 20x5 arithmetic practice\n// table (row-major order) with the new values from the commit diff.\n// The table cell text is fully replaced per-cell so existing run\n// formatting (TimeNewRoman, sz 30) on each paragraph is preserved.\nconst newValues = [\n  [\"69+25=94\", \"45-29=16\", \"90-88=2\", \"33-28=5\", \"23-6=17\"],\n  [\"23+29=52\", \"27+15=42\", \"35+49=84\", \"64-7=57\", \"85-6=79\"],\n  [\"72-39=33\", \"72-19=53\", \"18+27=45\", \"83-15=68\", \"72-49=23\"],\n  [\"60-48=12\", \"58-19=39\", \"97-69=28\", \"27+4=31\", \"71-58=13\"],\n  [\"10-3=7\", \"14+27=41\", \"90-74=16\", \"93-49=44\", \"38+34=72\"],\n  [\"71-68=3\", \"38+6=44\", \"15+17=32\", \"24-19=5\", \"95-38=57\"],\n  [\"34+8=42\", \"76-48=28\", \"14+28=42\", \"17+67=84\", \"81-7=74\"],\n  [\"70-1=69\", \"90-5=85\", \"33+48=81\", \"35-18=17\", \"80-59=21\"],\n  [\"2+19=21\", \"82-63=19\", \"29+44=73\", \"23+68=91\", \"23+29=52\"],\n  [\"6+79=85\", \"63-17=46\", \"39+2=41\", \"93-55=38\", \"83-48=35\"],\n  [\"71-68=3\", \"17+78=95\", \"59+8=67\", \"18+9=27\", \"33-17=16\"],\n  [\"11-3=8\", \"54-15=39\", \"2+49=51\", \"73-8=65\", \"5+16=21\"],\n  [\"27-18=9\", \"90-66=24\", \"67+6=73\", \"68+26=94\", \"47+7=54\"],\n  [\"56-28=28\", \"81-57=24\", \"74-5=69\", \"15+49=64\", \"58+34=92\"],\n  [\"81-25=56\", \"77+16=93\", \"83-59=24\", \"72-27=45\", \"90-83=7\"],\n  [\"48-39=9\", \"60-4=56\", \"31-5=26\", \"33+28=61\", \"84-28=56\"],\n  [\"75+19=94\", \"80-32=48\", \"85-16=69\", \"7+86=93\", \"23+19=42\"],\n  [\"14-7=7\", \"58+26=84\", \"80-56=24\", \"83-79=4\", \"42-36=6\"],\n  [\"32-4=28\", \"40-11=29\", \"28+14=42\", \"26+5=31\", \"62-53=9\"],\n  [\"27+59=86\", \"94-55=39\", \"45-29=16\", \"34+9=43\", \"8+37=45\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\n// Build the next grid by copying the current values and overwriting\n// each cell with its corresponding new value (keeps shape in sync\n// with the live table in case of any unexpected row/column count).\nconst currentValues = table.values;\nconst nextValues = currentValues.map((row, r) =>\n  row.map((cell, c) => (newValues[r] && newValues[r][c] !== undefined ? newValues[r][c] : cell))\n);\n\ntable.values = nextValues;\nawait context.sync();\n", "ps1": "# Replace the answer text in each cell of the 20x5 practice table\n# (row-major order), matching the new values from the commit diff.\n$newValues = @(\n  @(\"69+25=94\", \"45-29=16\", \"90-88=2\", \"33-28=5\", \"23-6=17\"),\n  @(\"23+29=52\", \"27+15=42\", \"35+49=84\", \"64-7=57\", \"85-6=79\"),\n  @(\"72-39=33\", \"72-19=53\", \"18+27=45\", \"83-15=68\", \"72-49=23\"),\n  @(\"60-48=12\", \"58-19=39\", \"97-69=28\", \"27+4=31\", \"71-58=13\"),\n  @(\"10-3=7\", \"14+27=41\", \"90-74=16\", \"93-49=44\", \"38+34=72\"),\n  @(\"71-68=3\", \"38+6=44\", \"15+17=32\", \"24-19=5\", \"95-38=57\"),\n  @(\"34+8=42\", \"76-48=28\", \"14+28=42\", \"17+67=84\", \"81-7=74\"),\n  @(\"70-1=69\", \"90-5=85\", \"33+48=81\", \"35-18=17\", \"80-59=21\"),\n  @(\"2+19=21\", \"82-63=19\", \"29+44=73\", \"23+68=91\", \"23+29=52\"),\n  @(\"6+79=85\", \"63-17=46\", \"39+2=41\", \"93-55=38\", \"83-48=35\"),\n  @(\"71-68=3\", \"17+78=95\", \"59+8=67\", \"18+9=27\", \"33-17=16\"),\n  @(\"11-3=8\", \"54-15=39\", \"2+49=51\", \"73-8=65\", \"5+16=21\"),\n  @(\"27-18=9\", \"90-66=24\", \"67+6=73\", \"68+26=94\", \"47+7=54\"),\n  @(\"56-28=28\", \"81-57=24\", \"74-5=69\", \"15+49=64\", \"58+34=92\"),\n  @(\"81-25=56\", \"77+16=93\", \"83-59=24\", \"72-27=45\", \"90-83=7\"),\n  @(\"48-39=9\", \"60-4=56\", \"31-5=26\", \"33+28=61\", \"84-28=56\"),\n  @(\"75+19=94\", \"80-32=48\", \"85-16=69\", \"7+86=93\", \"23+19=42\"),\n  @(\"14-7=7\", \"58+26=84\", \"80-56=24\", \"83-79=4\", \"42-36=6\"),\n  @(\"32-4=28\", \"40-11=29\", \"28+14=42\", \"26+5=31\", \"62-53=9\"),\n  @(\"27+59=86\", \"94-55=39\", \"45-29=16\", \"34+9=43\", \"8+37=45\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($r = 0; $r -lt $newValues.Length; $r++) {\n    $rowValues = $newValues[$r]\n    for ($c = 0; $c -lt $rowValues.Length; $c++) {\n        $t.Cell($r + 1, $c + 1).Range.Text = $rowValues[$c]\n    }\n}\n"}
